# Include the image title (from the markdown link title) in the
# PowerPoint description of the picture, alongside the link that was
# already present.  The picture "lalune.jpg" gains the "fig:  " prefix
# that pandoc uses to carry the title + link through to the pptx
# shape's alt-text / description field.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(1)
$shape.AlternativeText = "fig:  lalune.jpg"
